# Post-meeting updates to F2F presentation (2021-10-04 -> 2021-10-05)

$p = $ppt.ActivePresentation

# 1) Slides 2-6: "Date Placeholder 4" fields (2021-10-04 -> 2021-10-05)
for ($si = 2; $si -le 6; $si++) {
    $s = $p.Slides.Item($si)
    foreach ($shp in $s.Shapes) {
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "2021-10-05"
        }
    }
}

# 2) Slide 7: "Date Placeholder 5" field (2021-10-04 -> 2021-10-05)
$s7 = $p.Slides.Item(7)
foreach ($shp in $s7.Shapes) {
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = "2021-10-05"
    }
}

# 3) Slide master date placeholder (2021-10-04 -> 2021-10-05)
$master = $p.SlideMaster
foreach ($shp in $master.Shapes) {
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = "2021-10-05"
    }
}

# 4) All 11 slide layouts (layout 1 "Title Slide" has no date placeholder)
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $lay = $layouts.Item($li)
    foreach ($shp in $lay.Shapes) {
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "2021-10-05"
        }
    }
}

# 5) Notes master "datetimeFigureOut" field (10/4/21 -> 10/5/21)
$nm = $p.NotesMaster
$nmhf = $nm.HeadersFooters
$nmhf.DateAndTime.Text = "10/5/21"

# 6) Slide 7 content text edit: collapse last two bullets into one new bullet
$contentShape = $s7.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange
$para10 = $tr.Paragraphs(10, 1)
$para10.Text = "Need two implementations now of EVERY feature"
$tail = $para10.InsertAfter(", even optional ones")

$tr2 = $contentShape.TextFrame.TextRange
$para11 = $tr2.Paragraphs(11, 1)
$toRemove = $tr2.Characters($para11.Start, $para11.Length + 1)
$toRemove.Delete()
